# Add new chemical-database rows (Water ... Calcium nitrate) to the bottom
# of Sheet1, reusing the formatting already present in the last existing
# data row (row 50) so no superfluous cell-style records get created.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 51: the sheet already had an (empty) styled row here; only the
#     B:D cells picked up the "new data row" look, A51 kept its own style. ---
$ws.Range("A51").Value = "Water"
$ws.Range("B51").Value = "H₂O"
$ws.Range("C51").Value = 18.015
$ws.Range("D51").Value = 1

# --- Rows 52-65: brand-new rows, formatted like row 50. ---
$ws.Range("A52").Value = "Carbon dioxide"
$ws.Range("B52").Value = "CO₂"
$ws.Range("C52").Value = 44.01
$ws.Range("D52").Value = 1

# (formula typed in before the name here, matching the original authoring
#  order so the shared-string table indices line up)
$ws.Range("B53").Value = "LiCl"
$ws.Range("B54").Value = "LiBr"
$ws.Range("A53").Value = "Lithium chloride"
$ws.Range("A54").Value = "Lithium bromide"
$ws.Range("C53").Value = 42.39
$ws.Range("D53").Value = 1
$ws.Range("C54").Value = 86.85
$ws.Range("D54").Value = 1

$ws.Range("A55").Value = "Potassium chloride"
$ws.Range("B55").Value = "KCl"
$ws.Range("C55").Value = 74.55
$ws.Range("D55").Value = 1

$ws.Range("A56").Value = "Potassium bromide"
$ws.Range("B56").Value = "KBr"
$ws.Range("C56").Value = 119.002
$ws.Range("D56").Value = 1

$ws.Range("A57").Value = "Calcium chloride"
$ws.Range("B57").Value = "CaCl₂"
$ws.Range("C57").Value = 110.98
$ws.Range("D57").Value = 1

$ws.Range("A58").Value = "Calcium bromide"
$ws.Range("B58").Value = "CaBr₂"
$ws.Range("C58").Value = 199.89
$ws.Range("D58").Value = 1

$ws.Range("A59").Value = "Lithium sulfate"
$ws.Range("B59").Value = "Li₂SO₄"
$ws.Range("C59").Value = 109.94
$ws.Range("D59").Value = 1

$ws.Range("A60").Value = "Potassium sulfate"
$ws.Range("B60").Value = "K₂SO₄"
$ws.Range("C60").Value = 174.26
$ws.Range("D60").Value = 1

$ws.Range("A61").Value = "Calcium sulfate"
$ws.Range("B61").Value = "CaSO₄"
$ws.Range("C61").Value = 136.14
$ws.Range("D61").Value = 1

$ws.Range("A62").Value = "Nitrogen dioxide"
$ws.Range("B62").Value = "NO₂"
$ws.Range("C62").Value = 46.01
$ws.Range("D62").Value = 1

$ws.Range("A63").Value = "Lithium nitrate"
$ws.Range("B63").Value = "LiNO₃"
$ws.Range("C63").Value = 68.95
$ws.Range("D63").Value = 1

$ws.Range("A64").Value = "Potassium nitrate"
$ws.Range("B64").Value = "KNO₃"
$ws.Range("C64").Value = 101.1
$ws.Range("D64").Value = 1

$ws.Range("A65").Value = "Calcium nitrate"
$ws.Range("B65").Value = "Ca(NO₃)₂"
$ws.Range("C65").Value = 164.09
$ws.Range("D65").Value = 1

# --- Re-apply the existing row-50 formatting down the new block, the way
#     Excel does when you type a value straight below a formatted table:
#     column A (new rows only, row 51's A cell keeps its own old style),
#     and columns B:D (all new rows). This avoids inventing any new
#     cell-style (xf) records. ---
$ws.Range("A50").Copy() | Out-Null
$ws.Range("A52:A65").PasteSpecial(-4122) | Out-Null

$ws.Range("B50:D50").Copy() | Out-Null
$ws.Range("B51:D65").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- View state: the user scrolled down and zoomed in slightly before
#     saving, with the last edited cells (D62:D65) selected. ---
$ws.Range("D62:D65").Select() | Out-Null
$excel.ActiveWindow.Zoom = 102
